$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new columns before the old "most_frequent_value" column (Q),
# which shifts the existing Q,R,S,T columns (most_frequent_value,
# memory_consumed_bytes, pattern_count, patterns) two places right to
# S,T,U,V respectively.
$ws.Range("Q1:R1").EntireColumn.Insert()

# New header cells for the two inserted columns.
$ws.Range("Q1").Value = "default_count"
$ws.Range("R1").Value = "default_value"

# --- Row 2 (APPID) ---
$ws.Range("Q2").Value = 0
$ws.Range("R2").Value = "<Unspecified>"
$ws.Range("S2").Value = ""

# --- Row 3 (TypeOfViolence) ---
$ws.Range("Q3").Value = 0
$ws.Range("R3").Value = "<Unspecified>"
$ws.Range("S3").Value = "Assault"

# --- Row 4 (VictimType) ---
$ws.Range("Q4").Value = 0
$ws.Range("R4").Value = "<Unspecified>"
$ws.Range("S4").Value = "Primary"

# --- Row 5 (Gender) ---
$ws.Range("Q5").Value = 0
$ws.Range("R5").Value = "<Unspecified>"
$ws.Range("S5").Value = "F"

# --- Row 6 (Indigenous) ---
$ws.Range("Q6").Value = 0
$ws.Range("R6").Value = "<Unspecified>"
$ws.Range("S6").Value = "No"

# --- Row 7 (AgeAtHappenedTo) ---
$ws.Range("Q7").Value = 0
$ws.Range("R7").Value = "<Unspecified>"
# Stored as text (matches original convention where most_frequent_value
# keeps numeric-looking values as strings), so force a text format before
# assignment.
$ws.Range("S7").NumberFormat = "@"
$ws.Range("S7").Value = "0"

# --- Row 8 (Granted) ---
$ws.Range("Q8").Value = 3058
$ws.Range("R8").NumberFormat = "@"
$ws.Range("R8").Value = "0"
$ws.Range("S8").NumberFormat = "@"
$ws.Range("S8").Value = "0"
